$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Enable "Different Odd and Even Pages" headers/footers for this section.
# This mints the new "even" header/footer parts (rId8/rId10) and shifts the
# existing default/first header/footer parts to new relationship ids,
# matching Word's own renumbering behavior.
$sec.PageSetup.OddAndEvenPagesHeaderFooter = $true

# Touch the new even header/footer so Word actually materializes their
# parts (header1.xml / footer1.xml) as blank placeholders, and renumbers
# the remaining header/footer parts (header2/3.xml, footer2/3.xml).
$evenHeader = $sec.Headers(3)
$evenFooter = $sec.Footers(3)
$evenHeader.Range.Text = ""
$evenFooter.Range.Text = ""

# Update the "first page" header text (now header3.xml) with the teacher's
# name.
$firstHeader = $sec.Headers(2)
$firstHeader.Range.Find.Execute("Teacher Name", $true, $false, $false, $false, $false,
                                 $true, 1, $false, "Mrs. Patil", 2)
